$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 86, pushing existing rows 86.. down to 88..
$ws.Rows("86:87").Insert()

# New row 86 values
$ws.Range("A86").Value = 7
$ws.Range("B86").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C86").Value = "Ñuble"
$ws.Range("D86").Value = 44539
$ws.Range("E86").Value = 16
$ws.Range("F86").Value = "Fruta"
$ws.Range("G86").Value = 100101
$ws.Range("H86").Value = "Berries"
$ws.Range("I86").Value = 100112025
$ws.Range("J86").Value = "Frutilla"
$ws.Range("K86").Value = "Sin especificar"
$ws.Range("L86").Value = "Primera"
$ws.Range("M86").Value = 400
$ws.Range("N86").Value = 7000
$ws.Range("O86").Value = 7500
$ws.Range("P86").Value = 7250
$ws.Range("Q86").Value = "`$/caja 7 kilos"
$ws.Range("R86").Value = "Provincia de Diguillín"
$ws.Range("S86").Value = 1036
$ws.Range("T86").Value = 7

# New row 87 values
$ws.Range("A87").Value = 7
$ws.Range("B87").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C87").Value = "Ñuble"
$ws.Range("D87").Value = 44539
$ws.Range("E87").Value = 16
$ws.Range("F87").Value = "Fruta"
$ws.Range("G87").Value = 100101
$ws.Range("H87").Value = "Berries"
$ws.Range("I87").Value = 100112025
$ws.Range("J87").Value = "Frutilla"
$ws.Range("K87").Value = "Sin especificar"
$ws.Range("L87").Value = "Segunda"
$ws.Range("M87").Value = 200
$ws.Range("N87").Value = 6000
$ws.Range("O87").Value = 6500
$ws.Range("P87").Value = 6250
$ws.Range("Q87").Value = "`$/caja 7 kilos"
$ws.Range("R87").Value = "Provincia de Diguillín"
$ws.Range("S87").Value = 893
$ws.Range("T87").Value = 7
